$wb = $excel.ActiveWorkbook

# --- OFF sheet (offensive target depth data) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 215
$wsOff.Range("C2").Value = 151
$wsOff.Range("D2").Value = 46
$wsOff.Range("E2").Value = 25

# --- DEF sheet (defensive target depth data) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 171
$wsDef.Range("C2").Value = 119
$wsDef.Range("D2").Value = 40
$wsDef.Range("E2").Value = 26
$wsDef.Range("F2").Value = 3
$wsDef.Range("G2").Value = 2
